# Switch all labels and buttons to pics
# - Insert a new "gender" column before "number of exercises"
# - Populate gender values for existing patients
# - Update a few exercise flags that changed
# - Normalize patient ID in row 3 to a number
# - Update the active cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D ("gender"), shifting the existing
# "number of exercises" + exercise columns one to the right (E..U)
$ws.Columns("D:D").Insert()

# Header for the new column
$ws.Range("D1").Value = "gender"

# Gender values for the two existing patients
$ws.Range("D2").Value = "Male"
$ws.Range("D3").Value = "Male"

# A couple of exercise flags flipped on for patient in row 2
# (columns shifted right by one after the insert above)
$ws.Range("K2").Value = $true
$ws.Range("O2").Value = $true

# Patient ID in row 3 becomes a real number instead of text
$ws.Range("A3").Value = 11223344

# Move the active selection
$ws.Range("O8").Select() | Out-Null
